$d = $word.ActiveDocument

# The README's docker-compose snippets reference "01-slides/index.adoc" in
# four command blocks (asciidoctor-revealjs, asciidoctor-pdf,
# asciidoctor-epub3 and asciidoctor/docbook). Rename the referenced path to
# "_content/index.adoc" everywhere it appears. The surrounding double quotes
# are left untouched (kept out of the search/replace strings) so Word's
# smart-quote AutoCorrect doesn't turn them into curly quotes.
$d.Content.Find.Execute("01-slides/index.adoc", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "_content/index.adoc", 2)
